$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. C200 gets new value "已完成" (task for 黄宝怡 marked as finished)
$ws.Range("C200").Value = "已完成"

# 2. Row 204: new section header (merged A:D), styled like row 198
$ws.Range("A198:D198").Copy() | Out-Null
$ws.Range("A204:D204").PasteSpecial(-4122) | Out-Null
$ws.Range("A204").Value = "日期：2017.11.13-14 十二 周一  下午"
$ws.Rows.Item(204).RowHeight = 22.5
$ws.Range("A204:D204").Merge() | Out-Null

# 3. Row 205: column-titles row, styled like row 199
$ws.Range("A199:D199").Copy() | Out-Null
$ws.Range("A205:D205").PasteSpecial(-4122) | Out-Null
$ws.Range("A205").Value = "人员"
$ws.Range("B205").Value = "计划任务"
$ws.Range("C205").Value = "完成情况"
$ws.Range("D205").Value = "备注"
$ws.Rows.Item(205).RowHeight = 22.5

# 4. Row 206: new member task row for 魏仲凯, styled like row 200
$ws.Range("A200:D200").Copy() | Out-Null
$ws.Range("A206:D206").PasteSpecial(-4122) | Out-Null
$ws.Range("A206").Value = "魏仲凯"
$ws.Range("B206").Value = "负责我的收藏界面代码块"
$ws.Rows.Item(206).RowHeight = 22.5

# 5. Row 207: new member task row for 黄宝怡, styled like row 200
$ws.Range("A200:D200").Copy() | Out-Null
$ws.Range("A207:D207").PasteSpecial(-4122) | Out-Null
$ws.Range("A207").Value = "黄宝怡"
$ws.Range("B207").Value = "负责作品详细信息界面代码块"
$ws.Rows.Item(207).RowHeight = 22.5

# 6. Row 208: summary row (merged A:D), styled like row 201
$ws.Range("A201:D201").Copy() | Out-Null
$ws.Range("A208:D208").PasteSpecial(-4122) | Out-Null
$ws.Range("A208").Value = "总结："
$ws.Rows.Item(208).RowHeight = 22.5
$ws.Range("A208:D208").Merge() | Out-Null

# View state: select B206 (matches the workbook's saved selection)
$ws.Range("B206").Select() | Out-Null

Write-Output "done"
